# Generate Report for Handoff
#
# For the rows that were still "low" priority (not yet handed off) in both
# the zh-cn and de-de localization-status sheets, the handoff report run
# promoted them to "ht" (handed-off) and stamped a fresh handoff datetime.
#
#   zh-cn!E4:E7   "low" -> "ht"
#   zh-cn!H4:H7   "2016-08-21 08:38:31" -> "2016-08-21 08:38:50"
#
#   de-de!E4:E7   "low" -> "ht"
#   de-de!H4:H7   "2016-08-21 08:38:36" -> "2016-08-21 08:38:54"
#   Overview!G4:G7 "2016-08-21 08:38:36" -> "2016-08-21 08:38:54"
#   (Overview's "Latest HO Xliff Generate Date" shows the same timestamp
#    text as de-de's handoff datetime for these rows, so it is refreshed
#    alongside it)

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
for ($r = 4; $r -le 7; $r++) {
    $zhcn.Range("E$r").Value = "ht"
    $zhcn.Range("H$r").Value = "2016-08-21 08:38:50"
}

$dede = $wb.Worksheets.Item("de-de")
for ($r = 4; $r -le 7; $r++) {
    $dede.Range("E$r").Value = "ht"
    $dede.Range("H$r").Value = "2016-08-21 08:38:54"
}

# The de-de handoff timestamp above is the same text that the Overview
# sheet's "Latest HO Xliff Generate Date" column already shows for these
# rows (they shared one string table entry), so refresh it here too.
$overview = $wb.Worksheets.Item("Overview")
for ($r = 4; $r -le 7; $r++) {
    $overview.Range("G$r").Value = "2016-08-21 08:38:54"
}
